# Add 2022-Q4 data: insert a new "2022-Q4" sheet (with the quarterly fund
# holdings detail) right after "总计", and add a corresponding summary row
# to the "总计" sheet, pushing the older quarters down by one row/position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet, name it "2022-Q4" and move it so it sits
#    right after "总计" (i.e. right before "2022-Q3").
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"

$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$newSheet.Move($wsQ3)

# NOTE: sheet object references captured before a Move() can end up
# pointing at the wrong sheet afterwards -- always re-fetch by name once
# the tab order has changed.
$q4 = $wb.Worksheets.Item("2022-Q4")

# Copy header-row / index-column formatting from "2022-Q3" so the new
# sheet's look matches its siblings.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Range("B1:H1").Copy($q4.Range("B1:H1"))
$wsQ3.Range("A2:A6").Copy($q4.Range("A2:A6"))

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with the fund holdings data.
# ---------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold text that looks numeric (fund codes, and
# numbers formatted to a fixed number of decimals) -- format the cells as
# text first so Excel doesn't silently convert them to numbers and drop
# leading/trailing zeros.
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $q4.Range($col + "2:" + $col + "6").NumberFormat = "@"
}

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "003567"
$q4.Range("C2").Value = "华夏行业景气混合"
$q4.Range("D2").Value = "109.60"
$q4.Range("E2").Value = "93.65"
$q4.Range("F2").Value = "2.72"
$q4.Range("G2").Value = "2.9811"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "006348"
$q4.Range("C3").Value = "银华盛利混合A"
$q4.Range("D3").Value = "11.81"
$q4.Range("E3").Value = "86.13"
$q4.Range("F3").Value = "5.58"
$q4.Range("G3").Value = "0.6590"
$q4.Range("H3").Value = 1

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "180013"
$q4.Range("C4").Value = "银华领先策略混合"
$q4.Range("D4").Value = "5.00"
$q4.Range("E4").Value = "90.56"
$q4.Range("F4").Value = "3.01"
$q4.Range("G4").Value = "0.1505"
$q4.Range("H4").Value = 7

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "015684"
$q4.Range("C5").Value = "银华盛利混合C"
$q4.Range("D5").Value = "1.68"
$q4.Range("E5").Value = "86.13"
$q4.Range("F5").Value = "5.58"
$q4.Range("G5").Value = "0.0937"
$q4.Range("H5").Value = 1

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "001728"
$q4.Range("C6").Value = "银华战略新兴灵活配置定期开放混合"
$q4.Range("D6").Value = "1.31"
$q4.Range("E6").Value = "91.78"
$q4.Range("F6").Value = "5.83"
$q4.Range("G6").Value = "0.0764"
$q4.Range("H6").Value = 2

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: add a new first data row for
#    2022-Q4 and shift the existing quarter rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the styled index column (A) down into the new row 6, copying
# the formatting already used by rows 2-5.
$total.Range("A5").Copy($total.Range("A6"))

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 4.9

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 3.32

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 4.96

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 12
$total.Range("D3").Value = 3.46

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 3.96
